# Update LoL language sheet: add 'credits' row after 'title', and replace
# the 'play' row with a 'new_game' row (continue/new game refactor).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 4 (right after the "title" row, row 3), shifting
# everything below it down by one. This makes room for the new "credits" key.
$ws.Rows(4).Insert()

$ws.Range("A4").Value = "credits"
$ws.Range("B4").Value = "Written by: David Dionisio"

# The old "play" / "PLAY" row has now shifted down from row 6 to row 7.
# Replace it in place with the new "new_game" / "NEW GAME" entry.
# (Set B before A so the shared-string table records "NEW GAME" ahead of
# "new_game", matching the saved workbook's string order.)
$ws.Range("B7").Value = "NEW GAME"
$ws.Range("A7").Value = "new_game"

# Update the sheet selection to match the saved workbook state.
$ws.Range("A3").Select()
